$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fix typo in task name (row 4): "Documenteación" -> "Documentación"
$ws.Range("A4").Value = "Generación de Documentación"

# Simplify task name (row 5): "la documentación" -> "Documentación"
$ws.Range("A5").Value = "Revisión y mejora de Documentación"

# Update schedule for "Preparación de la presentación del TFG" (row 6):
# start date moves from 11-Jul-2017 (42936) to 16-Dec-2023 (45274),
# duration moves from 3 to 20 days. End date (D6) recalculates via the
# existing shared formula (=B6+C6).
$ws.Range("B6").Value = 45274
$ws.Range("C6").Value = 20

# Update overall project end date ("Fin proyecto") from 45291 to 45342
$ws.Range("C15").Value = 45342

# Keep the Gantt chart's date axis maximum in sync with the new project end date
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$valueAxis = $chart.Axes(2)
$valueAxis.MaximumScale = 45342

# Restore the active cell selection to C5
$ws.Range("C5").Select()
